# DatosRegistrarInformeVisitaVerificacion.xlsx - regression data refresh
# (automation regression + Reprogramacion Otro Pagare feature/definitions advance)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two registered "Cod cliente" values used by the automation
# regression data set.
$ws.Range("A2").Value = "22114387"

# A3 loses its inherited "General" format in the source workbook; force it
# back to Text (same as A2) before writing so the numeric-looking id is
# kept/stored as a string (matches A2's style) instead of being coerced to
# a number.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "24681769"

# Leave the cursor where the author left it when they last saved the file.
$ws.Range("F10").Select()
